# Update data values in result_data_KNN sheet (commit: "Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.741999999999999
$ws.Range("B6").Value = 6.212
$ws.Range("B7").Value = 6.031999999999999
$ws.Range("D7").Value = -7.835000000000001
$ws.Range("B8").Value = 6
$ws.Range("D11").Value = -7.38
$ws.Range("D12").Value = -7.071000000000001
$ws.Range("E12").Value = 17.967
$ws.Range("E13").Value = 16.468
$ws.Range("E14").Value = 17.176
$ws.Range("D15").Value = -8.411
$ws.Range("B16").Value = 4.935
$ws.Range("E16").Value = 16.409
$ws.Range("E19").Value = 16.54
$ws.Range("B20").Value = 8.204000000000001
$ws.Range("D20").Value = -7.957000000000001
$ws.Range("E20").Value = 16.475
$ws.Range("B21").Value = 9.236000000000001
$ws.Range("D21").Value = -7.98
$ws.Range("D22").Value = -7.720000000000001
$ws.Range("E22").Value = 16.602
$ws.Range("D23").Value = -7.869999999999999
$ws.Range("B28").Value = 6.171
$ws.Range("B29").Value = 5.101
$ws.Range("D29").Value = -6.832000000000001
$ws.Range("B30").Value = 6.021000000000001
$ws.Range("B32").Value = 6.645999999999999
$ws.Range("D34").Value = -8.066999999999998
$ws.Range("E36").Value = 16.527
$ws.Range("B40").Value = 9.179
$ws.Range("D42").Value = -8.18
$ws.Range("D43").Value = -8.241
$ws.Range("E43").Value = 17.019
$ws.Range("D44").Value = -8.081
$ws.Range("D45").Value = -7.525000000000001
$ws.Range("B46").Value = 5.801
$ws.Range("D46").Value = -7.920999999999998
$ws.Range("E46").Value = 16.797
$ws.Range("D50").Value = -8.230999999999998
$ws.Range("E50").Value = 16.605
$ws.Range("B51").Value = 4.862
$ws.Range("D51").Value = -8.337
$ws.Range("B52").Value = 6.042
$ws.Range("B57").Value = 5.022999999999999
$ws.Range("D57").Value = -7.924000000000001
$ws.Range("B59").Value = 5.231999999999999
$ws.Range("B62").Value = 5.229
$ws.Range("D65").Value = -7.784999999999999
$ws.Range("B66").Value = 5.819
$ws.Range("D66").Value = -7.347
$ws.Range("D67").Value = -6.83
$ws.Range("B73").Value = 6.703
$ws.Range("B74").Value = 8.943999999999999
$ws.Range("E76").Value = 16.602
$ws.Range("B77").Value = 5.677
$ws.Range("D79").Value = -7.923
$ws.Range("D84").Value = -8.300000000000001
$ws.Range("D87").Value = -8.096
$ws.Range("B92").Value = 5.513
$ws.Range("D92").Value = -6.528
$ws.Range("E95").Value = 17.483
$ws.Range("D97").Value = -8.104000000000001
$ws.Range("E97").Value = 16.779
$ws.Range("E99").Value = 16.429
$ws.Range("B100").Value = 6.031000000000001
